# Update market-price derived columns (H-N) across several worksheets
# Source: scheduled market-data refresh (re-pulled currentAveragePrice / Leve profit figures)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20021.098
$ws.Range("I32").Value = 3458.5173
$ws.Range("J32").Value = 260178.5
$ws.Range("K32").Value = 3458.5173
$ws.Range("L32").Value = 260178.5
$ws.Range("M32").Value = -3171.5173
$ws.Range("N32").Value = -260752.5
$ws.Range("H61").Value = 2751.7715
$ws.Range("I61").Value = 2195.963
$ws.Range("J61").Value = 4627.625
$ws.Range("K61").Value = 2195.963
$ws.Range("L61").Value = 4627.625
$ws.Range("M61").Value = -1983.963
$ws.Range("N61").Value = -5051.625
$ws.Range("H132").Value = 2549.2896
$ws.Range("I132").Value = 1907.5807
$ws.Range("J132").Value = 5391.143
$ws.Range("K132").Value = 5722.742099999999
$ws.Range("L132").Value = 16173.429
$ws.Range("M132").Value = -3192.742099999999
$ws.Range("N132").Value = -21233.429
$ws.Range("H136").Value = 2751.7715
$ws.Range("I136").Value = 2195.963
$ws.Range("J136").Value = 4627.625
$ws.Range("K136").Value = 6587.889000000001
$ws.Range("L136").Value = 13882.875
$ws.Range("M136").Value = -4037.889000000001
$ws.Range("N136").Value = -18982.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2545.081
$ws.Range("I134").Value = 1892
$ws.Range("J134").Value = 5919.3335
$ws.Range("K134").Value = 5676
$ws.Range("L134").Value = 17758.0005
$ws.Range("M134").Value = -3141
$ws.Range("N134").Value = -22828.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1180
$ws.Range("I16").Value = 1225
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1225
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -938
$ws.Range("N16").Value = -1574
$ws.Range("H58").Value = 2729.4783
$ws.Range("I58").Value = 1268.6923
$ws.Range("J58").Value = 4628.5
$ws.Range("K58").Value = 1268.6923
$ws.Range("L58").Value = 4628.5
$ws.Range("M58").Value = -1065.6923
$ws.Range("N58").Value = -5034.5
$ws.Range("H113").Value = 1180
$ws.Range("I113").Value = 1225
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1225
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 945
$ws.Range("N113").Value = -5340
$ws.Range("H134").Value = 4688.143
$ws.Range("I134").Value = 2070.2856
$ws.Range("J134").Value = 7306
$ws.Range("K134").Value = 6210.8568
$ws.Range("L134").Value = 21918
$ws.Range("M134").Value = -3675.8568
$ws.Range("N134").Value = -26988
$ws.Range("H136").Value = 2729.4783
$ws.Range("I136").Value = 1268.6923
$ws.Range("J136").Value = 4628.5
$ws.Range("K136").Value = 3806.0769
$ws.Range("L136").Value = 13885.5
$ws.Range("M136").Value = -1256.0769
$ws.Range("N136").Value = -18985.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 7250
$ws.Range("I63").Value = 5000
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 15000
$ws.Range("L63").Value = 24000
$ws.Range("M63").Value = -14251
$ws.Range("N63").Value = -25498
$ws.Range("H64").Value = 3602.8
$ws.Range("I64").Value = 1008
$ws.Range("J64").Value = 4714.857
$ws.Range("K64").Value = 3024
$ws.Range("L64").Value = 14144.571
$ws.Range("M64").Value = -2754
$ws.Range("N64").Value = -14684.571
$ws.Range("H66").Value = 7250
$ws.Range("I66").Value = 5000
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 45000
$ws.Range("L66").Value = 72000
$ws.Range("M66").Value = -41256
$ws.Range("N66").Value = -79488
$ws.Range("H67").Value = 3602.8
$ws.Range("I67").Value = 1008
$ws.Range("J67").Value = 4714.857
$ws.Range("K67").Value = 3024
$ws.Range("L67").Value = 14144.571
$ws.Range("M67").Value = -2088
$ws.Range("N67").Value = -16016.571
$ws.Range("H70").Value = 3875
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 13500
$ws.Range("M70").Value = -5685
$ws.Range("N70").Value = -14130
$ws.Range("H73").Value = 3875
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 13500
$ws.Range("M73").Value = -4908
$ws.Range("N73").Value = -15684
$ws.Range("H107").Value = 459.03845
$ws.Range("I107").Value = 486.5
$ws.Range("J107").Value = 397.25
$ws.Range("K107").Value = 1459.5
$ws.Range("L107").Value = 1191.75
$ws.Range("M107").Value = 460.5
$ws.Range("N107").Value = -5031.75
$ws.Range("H129").Value = 1559.5333
$ws.Range("I129").Value = 1584.4445
$ws.Range("J129").Value = 1522.1666
$ws.Range("K129").Value = 4753.333500000001
$ws.Range("L129").Value = 4566.4998
$ws.Range("M129").Value = 246.6664999999994
$ws.Range("N129").Value = -14566.4998
$ws.Range("H140").Value = 4954.242
$ws.Range("I140").Value = 6425.7896
$ws.Range("J140").Value = 2957.1428
$ws.Range("K140").Value = 19277.3688
$ws.Range("L140").Value = 8871.428400000001
$ws.Range("M140").Value = -14097.3688
$ws.Range("N140").Value = -19231.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 27925
$ws.Range("J93").Value = 27925
$ws.Range("L93").Value = 27925
$ws.Range("N93").Value = -31669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8818.166999999999
$ws.Range("I136").Value = 6242.2856
$ws.Range("J136").Value = 12424.4
$ws.Range("K136").Value = 18726.8568
$ws.Range("L136").Value = 37273.2
$ws.Range("M136").Value = -16176.8568
$ws.Range("N136").Value = -42373.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 520500
$ws.Range("I5").Value = 2000
$ws.Range("J5").Value = 693333.3
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 693333.3
$ws.Range("M5").Value = -1888
$ws.Range("N5").Value = -693557.3
$ws.Range("H6").Value = 277776.75
$ws.Range("I6").Value = 50000
$ws.Range("J6").Value = 353702.34
$ws.Range("K6").Value = 50000
$ws.Range("L6").Value = 353702.34
$ws.Range("M6").Value = -49885
$ws.Range("N6").Value = -353932.34
$ws.Range("H22").Value = 9166.666999999999
$ws.Range("J22").Value = 9166.666999999999
$ws.Range("L22").Value = 9166.666999999999
$ws.Range("N22").Value = -9752.666999999999
$ws.Range("H136").Value = 2761.7222
$ws.Range("I136").Value = 1912
$ws.Range("J136").Value = 3611.4443
$ws.Range("K136").Value = 5736
$ws.Range("L136").Value = 10834.3329
$ws.Range("M136").Value = -3186
$ws.Range("N136").Value = -15934.3329
